$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns E..T (16 columns) for rows 2..10
$data = @{}
$data[2]  = @(3, 1, 11.905615, 35.716845, 0.8197078149061106, 0.8197078149061106, 3, 1, 4.913477, 14.740431, 0.3201889893918886, 0.3201889893918885, 58.497965473355, 526.481689260195, 0.2624614168514208, 0.2624614168514207)
$data[3]  = @(3, 1, 11.905615, 35.716845, 0.8197078149061106, 0.8197078149061106, 3, 1, 3.864911333333334, 11.594734, 0.2518587252793198, 0.2518587252793198, 46.01414634380333, 414.12731709423, 0.2064505653637496, 0.2064505653637496)
$data[4]  = @(3, 1, 11.905615, 35.716845, 0.8197078149061106, 0.8197078149061106, 3, 1, 6.567164333333333, 19.701493, 0.4279522853287917, 0.4279522853287916, 78.18613019439833, 703.6751717495849, 0.3507958326909402, 0.3507958326909402)
$data[5]  = @(3, 1, 1.63134, 4.894019999999999, 0.1123186115768849, 0.1123186115768849, 3, 1, 4.913477, 14.740431, 0.3201889893918886, 0.3201889893918885, 8.015551569179999, 72.13996412262, 0.03596318273070285, 0.03596318273070284)
$data[6]  = @(3, 1, 1.63134, 4.894019999999999, 0.1123186115768849, 0.1123186115768849, 3, 1, 3.864911333333334, 11.594734, 0.2518587252793198, 0.2518587252793198, 6.30498445452, 56.74486009068, 0.02828842233689728, 0.02828842233689728)
$data[7]  = @(3, 1, 1.63134, 4.894019999999999, 0.1123186115768849, 0.1123186115768849, 3, 1, 6.567164333333333, 19.701493, 0.4279522853287917, 0.4279522853287916, 10.71327786354, 96.41950077185999, 0.04806700650928477, 0.04806700650928476)
$data[8]  = @(3, 1, 0.9872630000000001, 2.961789, 0.0679735735170045, 0.0679735735170045, 3, 1, 4.913477, 14.740431, 0.3201889893918886, 0.3201889893918885, 4.850894043451001, 43.65804639105901, 0.02176438980976491, 0.02176438980976491)
$data[9]  = @(3, 1, 0.9872630000000001, 2.961789, 0.0679735735170045, 0.0679735735170045, 3, 1, 3.864911333333334, 11.594734, 0.2518587252793198, 0.2518587252793198, 3.815683957680667, 34.34115561912601, 0.01711973757867288, 0.01711973757867288)
$data[10] = @(3, 1, 0.9872630000000001, 2.961789, 0.0679735735170045, 0.0679735735170045, 3, 1, 6.567164333333333, 19.701493, 0.4279522853287917, 0.4279522853287916, 6.483518361219668, 58.35166525097701, 0.02908944612856671, 0.02908944612856671)

$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($r in 2..10) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $ws.Range("$col$r").Value = $vals[$i]
    }
}
